$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93:101 down to 94:102
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly record
$ws.Range("A93").Value = 7
$ws.Range("B93").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C93").Value = "Ñuble"
$ws.Range("D93").Value = 45124
$ws.Range("E93").Value = 16
$ws.Range("F93").Value = 100112013
$ws.Range("G93").Value = "Alcachofa"
$ws.Range("H93").Value = "Argentina(o)"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 50
$ws.Range("K93").Value = 17000
$ws.Range("L93").Value = 17000
$ws.Range("M93").Value = 17000
$ws.Range("N93").Value = "$/caja 50 unidades"
$ws.Range("O93").Value = "Provincia de Limarí"
$ws.Range("P93").Value = 340
$ws.Range("Q93").Value = 50
$ws.Range("R93").Value = "Hortaliza"
